$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I2").Value = 150.85715
$ws.Range("J2").Value = 1439
$ws.Range("K2").Value = 150.85715
$ws.Range("L2").Value = 1439
$ws.Range("M2").Value = -37.85714999999999
$ws.Range("N2").Value = -1665
$ws.Range("H11").Value = 160.22223
$ws.Range("I11").Value = 160.22223
$ws.Range("K11").Value = 160.22223
$ws.Range("M11").Value = -20.22223
$ws.Range("H28").Value = 1395.2667
$ws.Range("I28").Value = 1563.1818
$ws.Range("J28").Value = 933.5
$ws.Range("K28").Value = 1563.1818
$ws.Range("L28").Value = 933.5
$ws.Range("M28").Value = -1078.1818
$ws.Range("N28").Value = -1903.5
$ws.Range("H38").Value = 1075
$ws.Range("I38").Value = 225.14285
$ws.Range("J38").Value = 4049.5
$ws.Range("K38").Value = 675.4285500000001
$ws.Range("L38").Value = 12148.5
$ws.Range("M38").Value = -303.4285500000001
$ws.Range("N38").Value = -12892.5
$ws.Range("H42").Value = 251.4375
$ws.Range("I42").Value = 40
$ws.Range("J42").Value = 462.875
$ws.Range("K42").Value = 120
$ws.Range("L42").Value = 1388.625
$ws.Range("M42").Value = 110
$ws.Range("N42").Value = -1848.625
$ws.Range("H62").Value = 6957.381
$ws.Range("J62").Value = 7047
$ws.Range("L62").Value = 7047
$ws.Range("N62").Value = -8295
$ws.Range("H65").Value = 6957.381
$ws.Range("J65").Value = 7047
$ws.Range("L65").Value = 35235
$ws.Range("N65").Value = -41475
$ws.Range("H116").Value = 2425
$ws.Range("H125").Value = 6000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 6000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 54000
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -58920
$ws.Range("H137").Value = 2347
$ws.Range("I137").Value = 2434.5386
$ws.Range("J137").Value = 2119.4
$ws.Range("K137").Value = 7303.6158
$ws.Range("L137").Value = 6358.200000000001
$ws.Range("M137").Value = -4753.6158
$ws.Range("N137").Value = -11458.2
$ws.Range("H138").Value = 197969.73
$ws.Range("I138").Value = 40869.32
$ws.Range("J138").Value = 324663.62
$ws.Range("K138").Value = 122607.96
$ws.Range("L138").Value = 973990.86
$ws.Range("M138").Value = -117467.96
$ws.Range("N138").Value = -984270.86

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 10000
$ws.Range("I28").Value = 10000
$ws.Range("K28").Value = 10000
$ws.Range("M28").Value = -9808
$ws.Range("H32").Value = 13505.167
$ws.Range("I32").Value = 13505.167
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 13505.167
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -13218.167
$ws.Range("N32").ClearContents()
$ws.Range("H45").Value = 2799.158
$ws.Range("I45").Value = 2403.2307
$ws.Range("J45").Value = 3657
$ws.Range("K45").Value = 2403.2307
$ws.Range("L45").Value = 3657
$ws.Range("M45").Value = -2026.2307
$ws.Range("N45").Value = -4411
$ws.Range("H61").Value = 4611.375
$ws.Range("I61").Value = 3920.4688
$ws.Range("J61").Value = 7375
$ws.Range("K61").Value = 3920.4688
$ws.Range("L61").Value = 7375
$ws.Range("M61").Value = -3708.4688
$ws.Range("N61").Value = -7799
$ws.Range("H99").Value = 10000
$ws.Range("I99").Value = 10000
$ws.Range("K99").Value = 10000
$ws.Range("M99").Value = -7005
$ws.Range("H119").Value = 51131.668
$ws.Range("J119").Value = 51131.668
$ws.Range("L119").Value = 51131.668
$ws.Range("N119").Value = -60807.668
$ws.Range("H123").Value = 61500
$ws.Range("J123").Value = 61500
$ws.Range("L123").Value = 61500
$ws.Range("N123").Value = -71300
$ws.Range("H132").Value = 2364
$ws.Range("I132").Value = 2284.6965
$ws.Range("K132").Value = 6854.0895
$ws.Range("M132").Value = -4324.0895
$ws.Range("H136").Value = 4611.375
$ws.Range("I136").Value = 3920.4688
$ws.Range("J136").Value = 7375
$ws.Range("K136").Value = 11761.4064
$ws.Range("L136").Value = 22125
$ws.Range("M136").Value = -9211.4064
$ws.Range("N136").Value = -27225

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1031.6666
$ws.Range("I94").Value = 734.75
$ws.Range("K94").Value = 734.75
$ws.Range("M94").Value = -283.75
$ws.Range("H134").Value = 5054.204
$ws.Range("I134").Value = 1800.6
$ws.Range("J134").Value = 8443.375
$ws.Range("K134").Value = 5401.799999999999
$ws.Range("L134").Value = 25330.125
$ws.Range("M134").Value = -2866.799999999999
$ws.Range("N134").Value = -30400.125

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2583.4
$ws.Range("I31").Value = 1386.625
$ws.Range("J31").Value = 4711
$ws.Range("K31").Value = 1386.625
$ws.Range("L31").Value = 4711
$ws.Range("M31").Value = -1091.625
$ws.Range("N31").Value = -5301
$ws.Range("H34").Value = 2583.4
$ws.Range("I34").Value = 1386.625
$ws.Range("J34").Value = 4711
$ws.Range("K34").Value = 1386.625
$ws.Range("L34").Value = 4711
$ws.Range("M34").Value = -1184.625
$ws.Range("N34").Value = -5115
$ws.Range("H58").Value = 1480.7241
$ws.Range("I58").Value = 1397.76
$ws.Range("J58").Value = 1999.25
$ws.Range("K58").Value = 1397.76
$ws.Range("L58").Value = 1999.25
$ws.Range("M58").Value = -1194.76
$ws.Range("N58").Value = -2405.25
$ws.Range("H99").Value = 5444.25
$ws.Range("I99").Value = 3518.5715
$ws.Range("J99").Value = 8140.2
$ws.Range("K99").Value = 3518.5715
$ws.Range("L99").Value = 8140.2
$ws.Range("M99").Value = -2020.5715
$ws.Range("N99").Value = -11136.2
$ws.Range("H126").Value = 5444.25
$ws.Range("I126").Value = 3518.5715
$ws.Range("J126").Value = 8140.2
$ws.Range("K126").Value = 10555.7145
$ws.Range("L126").Value = 24420.6
$ws.Range("M126").Value = -8085.7145
$ws.Range("N126").Value = -29360.6
$ws.Range("H136").Value = 1480.7241
$ws.Range("I136").Value = 1397.76
$ws.Range("J136").Value = 1999.25
$ws.Range("K136").Value = 4193.28
$ws.Range("L136").Value = 5997.75
$ws.Range("M136").Value = -1643.28
$ws.Range("N136").Value = -11097.75

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 106
$ws.Range("I2").Value = 69.75
$ws.Range("K2").Value = 418.5
$ws.Range("M2").Value = -305.5
$ws.Range("H5").Value = 1640.2858
$ws.Range("I5").Value = 610
$ws.Range("J5").Value = 3014
$ws.Range("K5").Value = 1830
$ws.Range("L5").Value = 9042
$ws.Range("M5").Value = -1718
$ws.Range("N5").Value = -9266
$ws.Range("H135").Value = 1640.2858
$ws.Range("I135").Value = 610
$ws.Range("J135").Value = 3014
$ws.Range("K135").Value = 5490
$ws.Range("L135").Value = 27126
$ws.Range("M135").Value = -2955
$ws.Range("N135").Value = -32196

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 25000
$ws.Range("J26").Value = 25000
$ws.Range("L26").Value = 25000
$ws.Range("N26").Value = -25560
$ws.Range("H50").Value = 25000
$ws.Range("J50").Value = 25000
$ws.Range("L50").Value = 25000
$ws.Range("N50").Value = -25996
$ws.Range("H80").Value = 3655.0688
$ws.Range("I80").Value = 3599.4
$ws.Range("J80").Value = 3778.7778
$ws.Range("K80").Value = 3599.4
$ws.Range("L80").Value = 3778.7778
$ws.Range("M80").Value = -2601.4
$ws.Range("N80").Value = -5774.7778
$ws.Range("H82").Value = 80000
$ws.Range("J82").Value = 80000
$ws.Range("L82").Value = 80000
$ws.Range("N82").Value = -80766
$ws.Range("H83").Value = 3655.0688
$ws.Range("I83").Value = 3599.4
$ws.Range("J83").Value = 3778.7778
$ws.Range("K83").Value = 17997
$ws.Range("L83").Value = 18893.889
$ws.Range("M83").Value = -13005
$ws.Range("N83").Value = -28877.889
$ws.Range("H85").Value = 80000
$ws.Range("J85").Value = 80000
$ws.Range("L85").Value = 80000
$ws.Range("N85").Value = -82652
$ws.Range("H94").Value = 67840
$ws.Range("J94").Value = 67840
$ws.Range("L94").Value = 67840
$ws.Range("N94").Value = -69192
$ws.Range("H99").Value = 40444.285
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 40444.285
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 40444.285
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -44936.285

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4793.875
$ws.Range("I40").Value = 5092.852
$ws.Range("J40").Value = 4172.923
$ws.Range("K40").Value = 5092.852
$ws.Range("L40").Value = 4172.923
$ws.Range("M40").Value = -4956.852
$ws.Range("N40").Value = -4444.923
$ws.Range("H55").Value = 256.5238
$ws.Range("J55").Value = 299.8
$ws.Range("L55").Value = 299.8
$ws.Range("N55").Value = -645.8
$ws.Range("H107").Value = 9023.333000000001
$ws.Range("I107").Value = 9023.333000000001
$ws.Range("K107").Value = 9023.333000000001
$ws.Range("M107").Value = -7103.333000000001
$ws.Range("H118").Value = 86666.336
$ws.Range("J118").Value = 86666.336
$ws.Range("L118").Value = 86666.336
$ws.Range("N118").Value = -89980.336

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 434.09525
$ws.Range("I107").Value = 367.4375
$ws.Range("J107").Value = 647.4
$ws.Range("K107").Value = 1102.3125
$ws.Range("L107").Value = 1942.2
$ws.Range("M107").Value = 817.6875
$ws.Range("N107").Value = -5782.2
$ws.Range("H109").Value = 79333
$ws.Range("J109").Value = 79333
$ws.Range("L109").Value = 79333
$ws.Range("N109").Value = -82107
$ws.Range("H116").Value = 72450
$ws.Range("J116").Value = 72450
$ws.Range("L116").Value = 72450
$ws.Range("N116").Value = -81628
$ws.Range("H119").Value = 91849
$ws.Range("J119").Value = 91849
$ws.Range("L119").Value = 91849
$ws.Range("N119").Value = -101525
$ws.Range("H132").Value = 2321.48
$ws.Range("I132").Value = 2489.7334
$ws.Range("J132").Value = 2069.1
$ws.Range("K132").Value = 7469.2002
$ws.Range("L132").Value = 6207.299999999999
$ws.Range("M132").Value = -4939.2002
$ws.Range("N132").Value = -11267.3
